# Update battery import-specification: add a new "Discharging current DC"
# datapoint row (B_I_DIS_DC) right after the existing "Charging current DC"
# (B_I_DC) row, shifting all subsequent datapoint rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 26 (pushes existing rows 26..59 down to 27..60).
$ws.Rows("26:26").Insert()

# Populate the new row with the new datapoint definition.
$ws.Range("A26").Value = "datapoints"
$ws.Range("B26").Value = "B_I_DIS_DC"
$ws.Range("C26").Value = "A"
$ws.Range("D26").Value = "Discharging current DC"
